$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.212.40"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").Value = "2.489.29"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.11"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.96"
$ws.Range("E6").Value = "  -5.59%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -2.71%  "
$ws.Range("D9").Value = "2.488.83"
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("E10").Value = "  -4.12%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.95"
$ws.Range("E13").Value = "  -4.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.02"
$ws.Range("E14").Value = "  -3.79%  "
$ws.Range("D15").Value = "2.930.34"
$ws.Range("E15").Value = "  -2.49%  "
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("D17").Value = "67.032.30"
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("D18").Value = "2.491.14"
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "360.98"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("E23").Value = "  -6.26%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.86"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.85"
$ws.Range("E26").Value = "  -5.56%  "
$ws.Range("E27").Value = "  -8.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.994"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").Value = "2.612.20"
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").Value = "0.0₃0938"
$ws.Range("E30").Value = "  -6.17%  "
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "504.63"
$ws.Range("E32").Value = "  -7.07%  "
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("E34").Value = "  -5.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.23"
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.21"
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("E39").Value = "  -3.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.57"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.73"
$ws.Range("E41").Value = "  -4.75%  "
$ws.Range("E42").Value = "  -5.22%  "
$ws.Range("E43").Value = "  -5.95%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -3.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.33"
$ws.Range("E46").Value = "  -2.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.08"
$ws.Range("E47").Value = "  -3.96%  "
$ws.Range("E48").Value = "  -4.59%  "
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("E50").Value = "  -5.26%  "
$ws.Range("E51").Value = "  -4.23%  "
